# ---------------------------------------------------------------------------
# Commit: Wed, Aug 05, 2020  6:06:17 PM
#
# 1) The table on slide 16 switches from the custom "Table_0" table style
#    ({93EDA059-F05A-4249-BA29-C5DF13E56B3E}, defined in tableStyles.xml) to
#    the built-in table style {62370606-47A1-45B1-AC3F-7FD37A51A87A}.
#
# 2) The deck's theme (color scheme) switches from the "Integral" palette to
#    the stock "Office Theme" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/
#    folHlink all change; fonts stay Arial/Arial so the font scheme is
#    unaffected).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style -----------------------------------------------------
$slide = $p.Slides.Item(16)
foreach ($shp in $slide.Shapes) {
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{62370606-47A1-45B1-AC3F-7FD37A51A87A}")
    }
}

# --- 2) Theme colors: Integral -> Office Theme ---------------------------
function Set-SchemeColor($colorScheme, [int]$index, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorScheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$cs = $p.SlideMaster.ColorScheme

# Slot order matches VBA's ColorScheme.Colors(1..12):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
Set-SchemeColor $cs 1  "000000"
Set-SchemeColor $cs 2  "FFFFFF"
Set-SchemeColor $cs 3  "44546A"
Set-SchemeColor $cs 4  "E7E6E6"
Set-SchemeColor $cs 5  "5B9BD5"
Set-SchemeColor $cs 6  "ED7D31"
Set-SchemeColor $cs 7  "A5A5A5"
Set-SchemeColor $cs 8  "FFC000"
Set-SchemeColor $cs 9  "4472C4"
Set-SchemeColor $cs 10 "70AD47"
Set-SchemeColor $cs 11 "0563C1"
Set-SchemeColor $cs 12 "954F72"
